# API-Animal-PDD-Dispatcher Config.xlsx update
# Re-purposes the generic REFramework Config.xlsx template into the
# project-specific config: new OutputReport setting, new mail-related
# constants (Animals/mailSubject/mailAttachmentName/attachmentDownload/
# mailAttachmentName2), a real OrchestratorQueueName value, and the
# Mail.Manager asset row, plus the related view/selection bookkeeping.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Settings"
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("A2").Value = "OrchestratorQueueName"
$settings.Range("B2").Value = "APIQueue"

$settings.Range("A3").Value = "OrchestratorQueueFolder"
$settings.Range("C3").Value = "Folder name. The value must match a folder defined in Orchestrator and queue specified as OrchestratorQueueName should be created in this folder. For classic folders leave the value field empty."
$settings.Rows.Item(3).RowHeight = 45

$settings.Rows.Item(5).RowHeight = 30

$settings.Range("A7").Value = "OutputReport"
$settings.Range("B7").Value = "C:\Users\IonutVarga\Documents\UiPath\API-Animal-PDD-Dispatcher\Data\Input"

# ---------------------------------------------------------------------
# Sheet "Constants"
# ---------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

$constants.Rows.Item(2).RowHeight = 30
$constants.Range("C2").Value = "Must be 0 if working with Orchestrator queues. If > 0, the robot will retry the same transaction which failed with a system exception. Must be an integer value."

$constants.Range("A3").Value = "MaxConsecutiveSystemExceptions"
$constants.Range("C3").Value = "The number of consecutive system exceptions allowed. If MaxConsecutiveSystemExceptions is reached, the job is stopped. To disable this feature, set the value to 0. "
$constants.Rows.Item(3).RowHeight = 45

$constants.Range("C7").Value = "Static part of logging message. Calling Get Transaction Data."

$constants.Range("C9").Value = "Static part of logging message. Processed Transaction succesful."
$constants.Range("C10").Value = "Static part of logging message. Processed Transaction failed with business exception."
$constants.Range("C11").Value = "Static part of logging message. Processed Transaction failed with application exception."

$constants.Range("A12").Value = "ExceptionMessage_ConsecutiveErrors"
$constants.Range("B12").Value = "The maximum number of consecutive system exceptions was reached. "
$constants.Range("C12").Value = "Error message in case MaxConsecutiveSystemExceptions number is reached."

$constants.Range("A18").Value = "Animals"
$constants.Range("B18").Value = "dog,cat"

$constants.Range("A19").Value = "mailSubject"
$constants.Range("B19").Value = "Poze pentru colegi"

$constants.Range("A20").Value = "mailAttachmentName"
$constants.Range("B20").Value = "Colegi.xlsx"

$constants.Range("A21").Value = "attachmentDownload"
$constants.Range("B21").Value = "C:\Users\IonutVarga\Documents\UiPath\API-Animal-PDD-Dispatcher\Data\Attachment"

$constants.Range("A22").Value = "mailAttachmentName2"
$constants.Range("B22").Value = "Colegi"

# ---------------------------------------------------------------------
# Sheet "Assets"
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

$assets.Range("C1").Value = "OrchestratorAssetFolder"
$assets.Range("A2").Value = "Mail.Manager"
$assets.Range("B2").Value = "APIdogs.mail.manager"

# ---------------------------------------------------------------------
# View / selection state (Constants becomes the active tab)
# ---------------------------------------------------------------------
$settings.Range("A7").Select()
$assets.Range("A2").Select()
$constants.Activate()
$constants.Range("B23").Select()
